$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apiTest")

# Make apiTest the active sheet / active tab (matches activeTab 0 -> 1,
# and tabSelected flipping from LoginTest to apiTest).
$ws.Activate()

# New header columns.
$ws.Range("E1").Value = "newFirstName"
$ws.Range("F1").Value = "newLastName"
$ws.Range("G1").Value = "newPass"

# New data row values for row 2.
$ws.Range("E2").Value = "Srdjan1"
$ws.Range("F2").Value = "Rados1"
$ws.Range("G2").Value = "Test123@"

# Row 1 grows to the same row height as the rest of the sheet.
$ws.Rows.Item(1).RowHeight = 13.8

# Adjust column widths for the newly added / resized columns.
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 20.666666666666668
$ws.Columns.Item(6).ColumnWidth = 19.166666666666668
$ws.Columns.Item(7).ColumnWidth = 15.999999999999998

# New active selection on apiTest.
$ws.Range("G2").Select() | Out-Null
